$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '61.118.48'
Set-TextCell 2 5 '  -3.83%  '

# Row 3
Set-TextCell 3 4 '2.976.46'
Set-TextCell 3 5 '  -3.47%  '

# Row 4
Set-TextCell 4 4 '0.999'
Set-TextCell 4 5 '  -0.12%  '

# Row 5
Set-TextCell 5 4 '543.81'

# Row 6
Set-TextCell 6 4 '130.20'
Set-TextCell 6 5 '  -7.03%  '

# Row 7
Set-TextCell 7 5 '  -0.02%  '

# Row 8
Set-TextCell 8 4 '2.970.68'
Set-TextCell 8 5 '  -3.48%  '

# Row 9
Set-TextCell 9 4 '0.491'
Set-TextCell 9 5 '  -2.98%  '

# Row 10
Set-TextCell 10 5 '  -7.78%  '

# Row 11
Set-TextCell 11 4 '5.88'
Set-TextCell 11 5 '  -7.44%  '

# Row 12
Set-TextCell 12 4 '0.440'
Set-TextCell 12 5 '  -3.97%  '

# Row 13
Set-TextCell 13 4 '0.0000217'
Set-TextCell 13 5 '  -3.29%  '

# Row 14
Set-TextCell 14 4 '33.34'
Set-TextCell 14 5 '  -4.92%  '

# Row 15
Set-TextCell 15 4 '3.458.83'
Set-TextCell 15 5 '  -3.51%  '

# Row 16
Set-TextCell 16 2 'TRON'
Set-TextCell 16 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 16 4 '0.109'
Set-TextCell 16 5 '  -2.86%  '

# Row 17
Set-TextCell 17 2 'WrappedBTC'
Set-TextCell 17 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 17 4 '61.104.54'
Set-TextCell 17 5 '  -3.85%  '

# Row 18
Set-TextCell 18 4 '2.972.76'

# Row 19
Set-TextCell 19 4 '6.52'
Set-TextCell 19 5 '  -2.21%  '

# Row 20
Set-TextCell 20 4 '469.91'
Set-TextCell 20 5 '  -1.24%  '

# Row 21
Set-TextCell 21 4 '12.95'
Set-TextCell 21 5 '  -3.93%  '

# Row 22
Set-TextCell 22 4 '0.656'
Set-TextCell 22 5 '  -6.56%  '

# Row 23
Set-TextCell 23 4 '6.89'
Set-TextCell 23 5 '  -3.08%  '

# Row 24
Set-TextCell 24 4 '79.24'
Set-TextCell 24 5 '  +0.43%  '

# Row 25
Set-TextCell 25 4 '11.83'
Set-TextCell 25 5 '  -3.83%  '

# Row 26
Set-TextCell 26 4 '0.999'
Set-TextCell 26 5 '  -0.05%  '

# Row 27
Set-TextCell 27 4 '2.68'
Set-TextCell 27 5 '  -1.69%  '

# Row 28
Set-TextCell 28 4 '7.52'
Set-TextCell 28 5 '  -5.79%  '

# Row 29
Set-TextCell 29 4 '0.998'
Set-TextCell 29 5 '  -0.08%  '

# Row 30
Set-TextCell 30 5 '  -2.37%  '

# Row 31
Set-TextCell 31 4 '25.25'
Set-TextCell 31 5 '  -3.79%  '

# Row 32
Set-TextCell 32 5 '  -3.93%  '

# Row 33
Set-TextCell 33 4 '2.27'
Set-TextCell 33 5 '  -2.82%  '

# Row 34
Set-TextCell 34 4 '5.39'
Set-TextCell 34 5 '  -0.77%  '

# Row 35
Set-TextCell 35 4 '54.32'
Set-TextCell 35 5 '  -6.66%  '

# Row 36
Set-TextCell 36 4 '5.79'
Set-TextCell 36 5 '  -3.70%  '

# Row 37
Set-TextCell 37 4 '442.17'
Set-TextCell 37 5 '  -10.46%  '

# Row 38
Set-TextCell 38 4 '3.109.65'
Set-TextCell 38 5 '  -4.17%  '

# Row 39
Set-TextCell 39 4 '0.0781'
Set-TextCell 39 5 '  -2.47%  '

# Row 40
Set-TextCell 40 4 '0.0373'
Set-TextCell 40 5 '  -7.75%  '

# Row 41
Set-TextCell 41 5 '  -3.19%  '

# Row 42
Set-TextCell 42 4 '7.98'
Set-TextCell 42 5 '  -1.95%  '

# Row 43
Set-TextCell 43 5 '  -0.04%  '

# Row 44
Set-TextCell 44 4 '2.27'
Set-TextCell 44 5 '  -14.18%  '

# Row 45
Set-TextCell 45 4 '25.07'
Set-TextCell 45 5 '  -2.08%  '

# Row 46
Set-TextCell 46 4 '0.237'
Set-TextCell 46 5 '  -6.84%  '

# Row 47
Set-TextCell 47 4 '0.107'
Set-TextCell 47 5 '  -2.92%  '

# Row 48
Set-TextCell 48 5 '  +9.80%  '

# Row 49
Set-TextCell 49 5 '  -6.81%  '

# Row 50
Set-TextCell 50 4 '113.29'
Set-TextCell 50 5 '  -9.04%  '

# Row 51
Set-TextCell 51 4 '0.0₃0473'
Set-TextCell 51 5 '  -10.69%  '

$wb.Save()